$wb = $excel.ActiveWorkbook

# Sheet "展览" - update F column (想去人数 / number of people interested)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 312
$ws1.Range("F3").Value = 50
$ws1.Range("F4").Value = 476
$ws1.Range("F5").Value = 4572
$ws1.Range("F6").Value = 352
$ws1.Range("F9").Value = 706
$ws1.Range("F10").Value = 192

# Sheet "全部类型" - update F column (想去人数 / number of people interested)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 312
$ws4.Range("F3").Value = 50
$ws4.Range("F4").Value = 476
$ws4.Range("F5").Value = 4572
$ws4.Range("F6").Value = 352
$ws4.Range("F9").Value = 706
$ws4.Range("F11").Value = 192
